$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Class/Series labels to the corrected spelling
$ws.Range("B2").Value = "Series A"
$ws.Range("B3").Value = "Series B"

# Move the active selection to B4 (better validation UX - land on the series column)
$ws.Range("B4").Select()
